$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "32.766882,34.967053"
$ws.Range("B8").Value = "32.511729,35.502029"
$ws.Range("B11").Value = "32.792761,34.995336"
$ws.Range("B18").Value = "32.980490,35.542420"
$ws.Range("B19").Value = "31.960770,34.876512"
$ws.Range("B21").Value = "32.986934,35.708518"
$ws.Range("B24").Value = "33.005860,35.094090"
$ws.Range("B25").Value = "32.601426,35.289751"
$ws.Range("B29").Value = "33.194459,35.572940"
$ws.Range("B30").Value = "31.961063,34.807761"
$ws.Range("B31").Value = "32.036425,34.842884"
$ws.Range("B32").Value = "31.977527,34.808252"
$ws.Range("B34").Value = "31.784215,35.117210"
$ws.Range("B35").Value = "32.058998,34.815227"
$ws.Range("B36").Value = "32.121447,34.803699"
$ws.Range("B37").Value = "32.692764,34.940222"
$ws.Range("B41").Value = "32.775683,34.967878"
$ws.Range("B42").Value = "32.199671,35.212911"
$ws.Range("B43").Value = "32.049533,34.764483"
$ws.Range("B51").Value = "32.792761,34.995336"
$ws.Range("B56").Value = "31.814560,34.779980"
$ws.Range("B57").Value = "31.248833,35.198232"
$ws.Range("B58").Value = "32.074578,34.805974"
$ws.Range("B61").Value = "32.590574,34.936472"
$ws.Range("B64").Value = "32.018460,34.748167"
$ws.Range("B69").Value = "31.822668,35.253867"
$ws.Range("B70").Value = "31.791658,34.651074"
$ws.Range("B75").Value = "31.916670,35.016670"
$ws.Range("B77").Value = "32.095980,34.774333"
$ws.Range("B83").Value = "31.248833,35.198232"
$ws.Range("B85").Value = "31.945204,34.878075"
$ws.Range("B88").Value = "32.098181,34.896471"
$ws.Range("B91").Value = "31.755957,34.989832"
$ws.Range("B93").Value = "31.750898,35.207819"
$ws.Range("B95").Value = "31.858601,35.215336"
$ws.Range("B96").Value = "32.177911,34.905656"
$ws.Range("B97").Value = "32.049272,34.798714"
$ws.Range("B100").Value = "31.226237,34.809557"
$ws.Range("B101").Value = "31.068012,35.007848"
$ws.Range("B102").Value = "31.238084,34.794545"
$ws.Range("B104").Value = "32.009918,34.739188"
$ws.Range("B105").Value = "31.863239,34.743120"
$ws.Range("B106").Value = "31.928344,34.878259"
$ws.Range("B107").Value = "32.174304,34.930966"
$ws.Range("B108").Value = "32.045852,34.752438"
$ws.Range("B109").Value = "31.314100,34.620250"
$ws.Range("B111").Value = "32.860863,35.099385"
$ws.Range("B112").Value = "32.471921,34.946694"
$ws.Range("B113").Value = "31.663407,34.599960"
$ws.Range("B114").Value = "31.977527,34.808252"
$ws.Range("B116").Value = "31.682230,34.745240"
$ws.Range("B119").Value = "31.670900,34.779750"
$ws.Range("B120").Value = "32.093309,34.885509"
$ws.Range("B121").Value = "31.246177,34.808709"
$ws.Range("B123").Value = "31.756796,34.988601"
$ws.Range("B124").Value = "32.053920,34.770991"
